# This script updates the "Phases" column (column AF) values on Sheet1
# to fix the body-map phase labels (see commit message: "body map fix fail").
# Each cell below had its phase label reassigned to the correct phase.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF2").Value = "Dichotic_and_AFACT"
$ws.Range("AF3").Value = "MAB_phase"
$ws.Range("AF4").Value = "MAB_and_Digit_after"
$ws.Range("AF5").Value = "dichotic_phase"
$ws.Range("AF7").Value = "Dichotic_and_AFACT"
$ws.Range("AF9").Value = "Dichotic_and_AFACT"
$ws.Range("AF10").Value = "Dichotic_and_AFACT"
$ws.Range("AF13").Value = "MAB_and_AFACT"
$ws.Range("AF14").Value = "MAB_phase"
$ws.Range("AF16").Value = "MAB_and_AFACT"
$ws.Range("AF17").Value = "dichotic_phase"
$ws.Range("AF18").Value = "MAB_phase"
$ws.Range("AF20").Value = "MAB_phase"
$ws.Range("AF22").Value = "MAB_and_AFACT"
$ws.Range("AF23").Value = "dichotic_phase"
$ws.Range("AF24").Value = "Digit_before_and_AFACT"
$ws.Range("AF25").Value = "MAB_and_Digit_after"
$ws.Range("AF26").Value = "MAB_and_AFACT"
$ws.Range("AF27").Value = "dichotic_phase"
$ws.Range("AF28").Value = "MAB_phase"
$ws.Range("AF30").Value = "dichotic_phase"
$ws.Range("AF31").Value = "Digit_before_and_AFACT"
$ws.Range("AF32").Value = "Dichotic_and_AFACT"
$ws.Range("AF35").Value = "MAB_and_Digit_after"
$ws.Range("AF36").Value = "MAB_and_AFACT"
$ws.Range("AF37").Value = "Dichotic_and_AFACT"
$ws.Range("AF38").Value = "Digit_before_and_AFACT"
$ws.Range("AF39").Value = "Dichotic_and_AFACT"
$ws.Range("AF40").Value = "Dichotic_and_AFACT"
$ws.Range("AF43").Value = "dichotic_phase"
$ws.Range("AF44").Value = "Dichotic_and_AFACT"
$ws.Range("AF45").Value = "Dichotic_and_AFACT"
$ws.Range("AF46").Value = "MAB_and_Digit_after"
$ws.Range("AF48").Value = "dichotic_phase"
$ws.Range("AF49").Value = "MAB_phase"
$ws.Range("AF50").Value = "MAB_and_AFACT"
$ws.Range("AF51").Value = "MAB_and_Digit_after"
$ws.Range("AF52").Value = "MAB_and_AFACT"
$ws.Range("AF54").Value = "dichotic_phase"
$ws.Range("AF55").Value = "MAB_phase"
$ws.Range("AF56").Value = "Digit_before_and_AFACT"
$ws.Range("AF57").Value = "MAB_and_AFACT"
$ws.Range("AF63").Value = "MAB_phase"
$ws.Range("AF64").Value = "Dichotic_and_AFACT"
$ws.Range("AF65").Value = "MAB_and_Digit_after"
$ws.Range("AF66").Value = "Dichotic_and_AFACT"
$ws.Range("AF67").Value = "Dichotic_and_AFACT"
$ws.Range("AF68").Value = "MAB_phase"
$ws.Range("AF69").Value = "Dichotic_and_AFACT"
$ws.Range("AF70").Value = "MAB_and_AFACT"
$ws.Range("AF72").Value = "dichotic_phase"
$ws.Range("AF73").Value = "dichotic_phase"
$ws.Range("AF74").Value = "Digit_before_and_AFACT"
$ws.Range("AF75").Value = "MAB_phase"
$ws.Range("AF76").Value = "MAB_and_AFACT"
$ws.Range("AF78").Value = "Dichotic_and_AFACT"
$ws.Range("AF79").Value = "MAB_and_Digit_after"
$ws.Range("AF80").Value = "dichotic_phase"
$ws.Range("AF81").Value = "dichotic_phase"
